$d = $word.ActiveDocument

# Paragraph 1 currently holds: "This is a Microsoft word document."
# Target: same paragraph, but with three additional runs appended:
#   " (", "Changed main", ")"
# We build this as 3 separate trailing w:r elements (not merged into the
# existing run, and not merged with each other) by temporarily splitting
# the new text across 3 new paragraphs, then re-joining those paragraphs
# back into paragraph 1 by deleting the paragraph marks between them.
# Deleting a paragraph mark merges text but keeps each side's runs intact,
# which is how we get 4 distinct <w:r> elements in the final paragraph.

$p1 = $d.Paragraphs(1)
$r = $p1.Range

# $r.End is one position past paragraph 1's own trailing paragraph mark
# (i.e. the offset where paragraph 2 begins). Use it as the anchor for all
# absolute offsets below.
$base = $r.End

# Create 3 fresh paragraphs right after paragraph 1 and fill each with one
# of the new text segments.
$r.InsertParagraphAfter()
$segA = $d.Paragraphs(2)
$segA.Range.InsertAfter(" (")

$segA.Range.InsertParagraphAfter()
$segB = $d.Paragraphs(3)
$segB.Range.InsertAfter("Changed main")

$segB.Range.InsertParagraphAfter()
$segC = $d.Paragraphs(4)
$segC.Range.InsertAfter(")")

# Layout right now (absolute character offsets from the start of the doc):
#   [0, base-1)          "This is a Microsoft word document."
#   base-1                paragraph mark #1 (ends paragraph 1)
#   [base, base+2)         " ("
#   base+2                 paragraph mark #2 (ends the " (" paragraph)
#   [base+3, base+15)      "Changed main"
#   base+15                paragraph mark #3 (ends the "Changed main" paragraph)
#   [base+16, base+17)     ")"
#   base+17                paragraph mark #4 (ends the ")" paragraph -- kept)
#
# Delete the 3 paragraph marks (working from the end backwards so earlier
# offsets stay valid) to splice everything back into a single paragraph
# while keeping each segment as its own run.

$mark3 = $d.Range($base + 15, $base + 16)
$mark3.Delete()

$mark2 = $d.Range($base + 2, $base + 3)
$mark2.Delete()

$mark1 = $d.Range($base - 1, $base)
$mark1.Delete()
